# Updated symbol list on Thu Jan 12 21:38:07 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto-coin rows that moved since the last snapshot. The sheet stores
# these as plain text (e.g. "287.28", "2.86%"), so each target cell is
# pre-formatted as Text before the new value is written - this keeps the
# value a literal string instead of letting Excel coerce the numeric-
# looking text into a real number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $NewValue)
    $cell = $Sheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
}

# Map of row -> @{ D = newPrice; E = newVolume }  (omit a key if that
# column did not change for the row)
$updates = [ordered]@{
    2  = @{ D = "287.28";      E = "2.86%" }
    3  = @{ D = "28.74";       E = "4.67%" }
    4  = @{ D = "5.048";       E = "4.48%" }
    5  = @{ D = "0.06662";     E = "4.35%" }
    6  = @{ D = "7.355";       E = "4.47%" }
    7  = @{ D = "3.390";       E = "1.91%" }
    8  = @{ D = "1.377";       E = "4.58%" }
    9  = @{                    E = "5.00%" }
    10 = @{ D = "0.1558";      E = "0.76%" }
    11 = @{ D = "0.06557";     E = "-0.13%" }
    12 = @{ D = "0.07640";     E = "1.42%" }
    13 = @{ D = "0.02951";     E = "0.28%" }
    14 = @{ D = "0.08979";     E = "-0.24%" }
    15 = @{ D = "0.001608";    E = "2.28%" }
    16 = @{ D = "0.04494";     E = "2.17%" }
    17 = @{ D = "0.0006462";   E = "0.19%" }
    18 = @{ D = "0.006358";    E = "5.13%" }
    19 = @{ D = "3.446" }
    20 = @{ D = "2.258";       E = "1.27%" }
    21 = @{ D = "0.3216";      E = "2.21%" }
    22 = @{ D = "0.1311";      E = "-2.95%" }
    23 = @{ D = "4.070";       E = "4.09%" }
    24 = @{ D = "0.1555";      E = "3.38%" }
    25 = @{ D = "0.001180";    E = "0.42%" }
    26 = @{ D = "0.004492";    E = "5.02%" }
    27 = @{ D = "0.0001249";   E = "5.89%" }
    28 = @{                    E = "-2.03%" }
    40 = @{ D = "0.04205" }
    41 = @{ D = "0.006717";    E = "1.14%" }
    42 = @{ D = "0.1254";      E = "-11.13%" }
    43 = @{ D = "0.002019";    E = "-3.39%" }
    44 = @{ D = "0.01231";     E = "11.83%" }
    45 = @{ D = "0.00005733";  E = "3.18%" }
    46 = @{                    E = "20.74%" }
    47 = @{ D = "0.01309";     E = "-29.24%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        Set-TextValue $ws "D$row" $cols["D"]
    }
    if ($cols.ContainsKey("E")) {
        Set-TextValue $ws "E$row" $cols["E"]
    }
}
